$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

$oldText = "Administrator, Miss Dina Nasr"
$newText = "Miss Dina Nasr, Administrator"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldText) {
        $cell.Value = $newText
    }
}
